# Updates crypto price/volume figures per the latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.261.65"
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("D3").Value = "2.608.03"
$ws.Range("E3").Value = "  +3.55%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'307.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("D6").Value = "'99.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.84%  "
$ws.Range("D7").Value = "'0.601"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.93%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.579"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.15%  "
$ws.Range("D10").Value = "'39.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.13%  "
$ws.Range("D11").Value = "'0.0843"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.30%  "
$ws.Range("D12").Value = "'54.17"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "'8.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.86%  "
$ws.Range("D14").Value = "3.009.10"
$ws.Range("E14").Value = "  +3.28%  "
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").Value = "2.613.18"
$ws.Range("E16").Value = "  +3.46%  "
$ws.Range("D17").Value = "'0.916"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.68%  "
$ws.Range("D18").Value = "'14.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("D19").Value = "46.395.62"
$ws.Range("E19").Value = "  +0.86%  "
$ws.Range("E20").Value = "  +2.59%  "
$ws.Range("D21").Value = "'12.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.12%  "
$ws.Range("D22").Value = "'6.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.36%  "
$ws.Range("D23").Value = "'71.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.75%  "
$ws.Range("D24").Value = "'272.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.25%  "
$ws.Range("E25").Value = "  +3.69%  "
$ws.Range("E26").Value = "  +4.04%  "
$ws.Range("D27").Value = "'29.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +23.24%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("E29").Value = "  -1.03%  "
$ws.Range("E30").Value = "  +3.08%  "
$ws.Range("E31").Value = "  +0.96%  "
$ws.Range("D32").Value = "'38.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.17%  "
$ws.Range("E33").Value = "  +6.88%  "
$ws.Range("E34").Value = "  -4.11%  "
$ws.Range("E35").Value = "  -0.83%  "
$ws.Range("E36").Value = "  +3.12%  "
$ws.Range("E37").Value = "  +1.88%  "
$ws.Range("D38").Value = "'150.67"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.00%  "
$ws.Range("E39").Value = "  +4.49%  "
$ws.Range("D40").Value = "'0.123"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.10%  "
$ws.Range("D41").Value = "'23.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +32.15%  "
$ws.Range("D42").Value = "'15.93"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.26%  "
$ws.Range("D43").Value = "'3.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.10%  "
$ws.Range("D44").Value = "'0.0330"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.91%  "
$ws.Range("D45").Value = "'4.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.24%  "
$ws.Range("D46").Value = "2.113.86"
$ws.Range("E46").Value = "  +5.88%  "
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("D48").Value = "'93.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("D49").Value = "'9.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.79%  "
$ws.Range("E50").Value = "  -2.91%  "
$ws.Range("D51").Value = "'108.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.86%  "
